$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header "Pellet Count" (F1): right-align with an explicit (re-applied)
# Calibri 11 font, matching the column's now-numeric data below. ---
$f1 = $ws.Range("F1")
$f1.Font.Name = "Calibri"
$f1.Font.Size = 11
$f1.Borders.LineStyle = -4142   # xlLineStyleNone (no visible change, explicit touch)
$f1.HorizontalAlignment = -4152 # xlRight

# --- Attendance data rows 2-7 ---
# Columns B/C/D/E hold values that look numeric/date-like but must stay text,
# so they are entered with a leading apostrophe (the normal Excel way of
# forcing text for digit-only / date-like entries). Column F holds real
# numbers, and row 7's Time In / Time Out were entered as real numbers too.

function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

# Row 2
$ws.Range("A2").Value = "luth"
Set-TextValue $ws.Range("B2") "11/9/2022"
Set-TextValue $ws.Range("C2") "1600"
Set-TextValue $ws.Range("D2") "1800"
Set-TextValue $ws.Range("E2") "745305"
$ws.Range("F2").Value = 100

# Row 3
$ws.Range("A3").Value = "test"
Set-TextValue $ws.Range("B3") "11/9/2022"
Set-TextValue $ws.Range("C3") "1234"
Set-TextValue $ws.Range("D3") "1345"
Set-TextValue $ws.Range("E3") "123456"
$ws.Range("F3").Value = 100

# Row 4
$ws.Range("A4").Value = "luth"
Set-TextValue $ws.Range("B4") "11/9/2022"
Set-TextValue $ws.Range("C4") "1600"
Set-TextValue $ws.Range("D4") "1800"
Set-TextValue $ws.Range("E4") "745305"
$ws.Range("F4").Value = 100

# Row 5
$ws.Range("A5").Value = "test"
Set-TextValue $ws.Range("B5") "11/9/2022"
Set-TextValue $ws.Range("C5") "1234"
Set-TextValue $ws.Range("D5") "1345"
Set-TextValue $ws.Range("E5") "123456"
$ws.Range("F5").Value = 100

# Row 6
$ws.Range("A6").Value = "luth"
Set-TextValue $ws.Range("B6") "17/9/2022"
Set-TextValue $ws.Range("C6") "1600"
Set-TextValue $ws.Range("D6") "1800"
Set-TextValue $ws.Range("E6") "123456"
$ws.Range("F6").Value = 123

# Row 7
$ws.Range("A7").Value = "test"
Set-TextValue $ws.Range("B7") "17/9/2022"
$ws.Range("C7").Value = 1600
$ws.Range("D7").Value = 1800
Set-TextValue $ws.Range("E7") "123456"
$ws.Range("F7").Value = 123
